$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)

# Resize / reposition the content placeholder (it had an empty <p:spPr/> before).
$sh.Left = 30
$sh.Top = 93
$sh.Width = 654
$sh.Height = 354

$tr = $sh.TextFrame.TextRange
$cr = [char]13

# --- Bump the top-level bullet font size from 20pt to 18pt (first bullet). ---
$para1 = $tr.Paragraphs(1, 1)
$para1.Font.Size = 18

# --- Insert two new sub-bullets right after bullet 1, before "Once connected...". ---
$subText1 = "The ScienceBase folder ID can be found by accessing the folder you want on the ScienceBase website and looking at the last part of the URL"
$ins1 = $para1.InsertAfter($cr + $subText1)

$para2 = $tr.Paragraphs(2, 1)
$para2.IndentLevel = 2
$para2.Font.Size = 16
# Mark the two "ScienceBase" occurrences inside this paragraph as flagged words to mirror the source formatting.
$sb1Start = $para2.Start + ("The ").Length
$tr.Characters($sb1Start, ("ScienceBase").Length).Text = "ScienceBase"

$subText2 = "For example, for the URL: https://www.sciencebase.gov/catalog/item/5fbe75fad34e4b9faad7e8a1 the ScienceBase folder ID is " + [char]8220 + "5fbe75fad34e4b9faad7e8a1" + [char]8221 + " without the quotes."
$ins2 = $para2.InsertAfter($cr + $subText2)

$para3 = $tr.Paragraphs(3, 1)
$para3.IndentLevel = 3
$para3.Font.Size = 14

# Give the "the URL" run its own (non-dirty) formatting boundary, matching the source run split.
$urlWordStart = $para3.Start + ("For example, for ").Length
$urlWordLen = ("the URL").Length
$urlWordRange = $tr.Characters($urlWordStart, $urlWordLen)
$urlWordRange.Font.Size = 14

# Apply the hyperlink + color override to the literal URL text.
$urlStart = $para3.Start + ("For example, for the URL: ").Length
$urlLen = ("https://www.sciencebase.gov/catalog/item/5fbe75fad34e4b9faad7e8a1").Length
$urlRange = $tr.Characters($urlStart, $urlLen)
$urlRange.Font.Size = 14
$urlRange.ActionSettings(1).Hyperlink.Address = "https://www.sciencebase.gov/catalog/item/5fbe75fad34e4b9faad7e8a1"

# The single trailing space after the link keeps the same color override as the link run.
$spAfterStart = $urlStart + $urlLen
$spAfterRange = $tr.Characters($spAfterStart, 1)
$spAfterRange.Font.Size = 14

# --- Paragraph 4 is the unchanged "Once connected..." bullet; leave as-is. ---

# --- Remaining top-level bullets: bump 20pt -> 18pt. ---
for ($i = 5; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.IndentLevel -eq 1 -and $para.Font.Size -eq 20) {
        $para.Font.Size = 18
    }
}
